$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "trainingaudio/03_kikita3.wav"
$ws.Range("B2").Value = "pngimages/03_box.png"

$ws.Range("A3").Value = "trainingaudio/11_tokiko1.wav"
$ws.Range("B3").Value = "pngimages/11_compass.png"

$ws.Range("A4").Value = "trainingaudio/25_tapapi1.wav"
$ws.Range("B4").Value = "pngimages/25_apple.png"

$ws.Range("A5").Value = "trainingaudio/18_popata2.wav"
$ws.Range("B5").Value = "pngimages/18_donut.png"

$ws.Range("A6").Value = "trainingaudio/16_kotapi2.wav"
$ws.Range("B6").Value = "pngimages/16_icecream.png"

$ws.Range("A7").Value = "trainingaudio/10_tokiti1.wav"
$ws.Range("B7").Value = "pngimages/10_backpack.png"
